$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names for rows 2..26 (column A)
$names = @(
    "model_2_7_0",
    "model_2_7_22",
    "model_2_7_21",
    "model_2_7_20",
    "model_2_7_19",
    "model_2_7_18",
    "model_2_7_17",
    "model_2_7_16",
    "model_2_7_15",
    "model_2_7_14",
    "model_2_7_13",
    "model_2_7_23",
    "model_2_7_12",
    "model_2_7_10",
    "model_2_7_9",
    "model_2_7_8",
    "model_2_7_7",
    "model_2_7_6",
    "model_2_7_5",
    "model_2_7_4",
    "model_2_7_3",
    "model_2_7_2",
    "model_2_7_1",
    "model_2_7_11",
    "model_2_7_24"
)

# Common metric values shared by every row (B..I)
$values = @(
    0.09932080507882668,
    0.1341933169903529,
    0.4512388185747046,
    0.2880223605302709,
    0.9967864155769348,
    1.593163371086121,
    0.3132980167865753,
    0.9908738732337952
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
